$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J12: average of column J (k value), bold font (matches existing bold header font)
$j12 = $ws.Range("J12")
$j12.Formula = "=AVERAGE(J2:J11)"
$j12.Font.Bold = $true

# Row 14: Average of SW(S*)/SW(OPT)
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$b14 = $ws.Range("B14")
$b14.Formula = "=AVERAGE(N2:N11)"
$b14.Font.Bold = $true
$b14.Font.Size = 12
$b14.VerticalAlignment = -4108

# Row 15: Average of SC(S*)/SC(OPT)
$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$b15 = $ws.Range("B15")
$b15.Formula = "=AVERAGE(Z2:Z11)"
$b15.Font.Bold = $true
$b15.Font.Size = 12
$b15.VerticalAlignment = -4108

# Row 16: Worst of SW(S*)/SW(OPT)
$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$b16 = $ws.Range("B16")
$b16.Formula = "=MIN(N2:N11)"
$b16.Font.Bold = $true
$b16.Font.Size = 12
$b16.VerticalAlignment = -4108

# Row 17: Worst of SC(S*)/SC(OPT)
$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$b17 = $ws.Range("B17")
$b17.Formula = "=MAX(Z2:Z11)"
$b17.Font.Bold = $true
$b17.Font.Size = 12
$b17.VerticalAlignment = -4108

Write-Host "done"
